$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 37
$ws.Range("A37").Value = "2024-04-07 00:23:03"
$ws.Range("B37").Value = 4
$ws.Range("C37").Value = 4
$ws.Range("D37").Value = 0
$ws.Range("E37").Value = 2
$ws.Range("F37").Value = 2
$ws.Range("G37").Value = 0
$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0.001
$ws.Range("J37").Value = 0.05
$ws.Range("K37").Value = 0.003
$ws.Range("L37").Value = 100
$ws.Range("M37").Value = 500
$ws.Range("N37").Value = 10
$ws.Range("O37").Value = 9
$ws.Range("P37").Value = 3
$ws.Range("Q37").Value = 1000
$ws.Range("R37").Value = 5
$ws.Range("S37").Value = 1
$ws.Range("T37").Value = 100
$ws.Range("U37").Value = 1
$ws.Range("V37").Value = "Data/bombay1.xlsx"
$ws.Range("W37").Value = 196000
$ws.Range("X37").Value = "No es Simulación"

# Row 38
$ws.Range("A38").Value = "2024-04-14 19:49:46"
$ws.Range("B38").Value = 2
$ws.Range("C38").Value = 2
$ws.Range("D38").Value = 0
$ws.Range("E38").Value = 1
$ws.Range("F38").Value = 0
$ws.Range("G38").Value = 1
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0.001
$ws.Range("J38").Value = 0.05
$ws.Range("K38").Value = 0.003
$ws.Range("L38").Value = 100
$ws.Range("M38").Value = 500
$ws.Range("N38").Value = 10
$ws.Range("O38").Value = 9
$ws.Range("P38").Value = 3
$ws.Range("Q38").Value = 1000
$ws.Range("R38").Value = 5
$ws.Range("S38").Value = 1
$ws.Range("T38").Value = 100
$ws.Range("U38").Value = 1
$ws.Range("V38").Value = "Data/bombay1.xlsx"
$ws.Range("W38").Value = 136000
$ws.Range("X38").Value = "No es Simulación"

# Row 39
$ws.Range("A39").Value = "2024-04-18 21:37:03"
$ws.Range("B39").Value = 23
$ws.Range("C39").Value = 16
$ws.Range("D39").Value = 2
$ws.Range("E39").Value = 4
$ws.Range("F39").Value = 3
$ws.Range("G39").Value = 7
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0.001
$ws.Range("J39").Value = 0.05
$ws.Range("K39").Value = 0.003
$ws.Range("L39").Value = 100
$ws.Range("M39").Value = 500
$ws.Range("N39").Value = 10
$ws.Range("O39").Value = 9
$ws.Range("P39").Value = 3
$ws.Range("Q39").Value = 500
$ws.Range("R39").Value = 5
$ws.Range("S39").Value = 1
$ws.Range("T39").Value = 90
$ws.Range("U39").Value = 0.6956521739130435
$ws.Range("V39").Value = "Data/bombay1.xlsx"
$ws.Range("W39").Value = -314500
$ws.Range("X39").Value = "No es Simulación"

# Row 40
$ws.Range("A40").Value = "2024-04-20 22:25:44"
$ws.Range("B40").Value = 188
$ws.Range("C40").Value = 34
$ws.Range("D40").Value = 4
$ws.Range("E40").Value = 13
$ws.Range("F40").Value = 17
$ws.Range("G40").Value = 0
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0.001
$ws.Range("J40").Value = 0.05
$ws.Range("K40").Value = 0.003
$ws.Range("L40").Value = 100
$ws.Range("M40").Value = 500
$ws.Range("N40").Value = 10
$ws.Range("O40").Value = 9
$ws.Range("P40").Value = 2
$ws.Range("Q40").Value = 500
$ws.Range("R40").Value = 1
$ws.Range("S40").Value = 0
$ws.Range("T40").Value = 10
$ws.Range("U40").Value = 0.1808510638297872
$ws.Range("V40").Value = "Data/bombay1.xlsx"
$ws.Range("W40").Value = 59500
$ws.Range("X40").Value = "No es Simulación"

# Row 41
$ws.Range("A41").Value = "2024-04-20 23:11:48"
$ws.Range("B41").Value = 67
$ws.Range("C41").Value = 8
$ws.Range("D41").Value = 3
$ws.Range("E41").Value = 3
$ws.Range("F41").Value = 2
$ws.Range("G41").Value = 0
$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0.001
$ws.Range("J41").Value = 0.05
$ws.Range("K41").Value = 0.003
$ws.Range("L41").Value = 100
$ws.Range("M41").Value = 500
$ws.Range("N41").Value = 10
$ws.Range("O41").Value = 9
$ws.Range("P41").Value = 2
$ws.Range("Q41").Value = 500
$ws.Range("R41").Value = 1
$ws.Range("S41").Value = 0
$ws.Range("T41").Value = 50
$ws.Range("U41").Value = 0.1194029850746269
$ws.Range("V41").Value = "Data/bombay1.xlsx"
$ws.Range("W41").Value = -41000
$ws.Range("X41").Value = "No es Simulación"

# Row 42
$ws.Range("A42").Value = "2024-04-21 01:04:36"
$ws.Range("B42").Value = 14
$ws.Range("C42").Value = 10
$ws.Range("D42").Value = 2
$ws.Range("E42").Value = 6
$ws.Range("F42").Value = 2
$ws.Range("G42").Value = 0
$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 0.001
$ws.Range("J42").Value = 0.05
$ws.Range("K42").Value = 0.003
$ws.Range("L42").Value = 100
$ws.Range("M42").Value = 500
$ws.Range("N42").Value = 10
$ws.Range("O42").Value = 9
$ws.Range("P42").Value = 2
$ws.Range("Q42").Value = 1000
$ws.Range("R42").Value = 3
$ws.Range("S42").Value = 1
$ws.Range("T42").Value = 90
$ws.Range("U42").Value = 0.7142857142857143
$ws.Range("V42").Value = "Data/bombay1.xlsx"
$ws.Range("W42").Value = 533000
$ws.Range("X42").Value = "No es Simulación"
